$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run of the missing-data sampler: two rows worth of samples are dropped
# from the bottom of the table (table now ends at row 33 instead of row 35),
# and the id/value columns for rows 26:33 plus the blanked-out cell in row 3
# are refreshed with the new draw.

# Drop the now-unused trailing rows so the table ends at row 33.
$ws.Range("A34:F35").EntireRow.Delete()

# Row 3
$ws.Cells.Item(3, 5).Value = "'"
$ws.Cells.Item(3, 5).Style = "Normal"
# Row 26
$ws.Cells.Item(26, 1).Value = "SC 5"
$ws.Cells.Item(26, 2).Value = -20.2
$ws.Cells.Item(26, 3).Value = 10.8
$ws.Cells.Item(26, 4).Value = -13.8
$ws.Cells.Item(26, 5).Value = -5
$ws.Cells.Item(26, 6).Value = 17.38
# Row 27
$ws.Cells.Item(27, 1).Value = "SC 101"
$ws.Cells.Item(27, 2).Value = -20.4
$ws.Cells.Item(27, 4).Value = -14.6
$ws.Cells.Item(27, 5).Value = -10
$ws.Cells.Item(27, 6).Value = 17
# Row 28
$ws.Cells.Item(28, 1).Value = "SC 105"
$ws.Cells.Item(28, 2).Value = -19.6
$ws.Cells.Item(28, 3).Value = 11.1
$ws.Cells.Item(28, 4).Value = -13.7
$ws.Cells.Item(28, 5).Value = -5.9
$ws.Cells.Item(28, 6).Value = 17.44
# Row 29
$ws.Cells.Item(29, 1).Value = "SC 119"
$ws.Cells.Item(29, 2).Value = -19.5
$ws.Cells.Item(29, 3).Value = 11.2
$ws.Cells.Item(29, 4).Value = -13
$ws.Cells.Item(29, 5).Value = -6.8
$ws.Cells.Item(29, 6).Value = 18.06
# Row 30
$ws.Cells.Item(30, 1).Value = "SC 120"
$ws.Cells.Item(30, 2).Value = -19.7
$ws.Cells.Item(30, 3).Value = 11.4
$ws.Cells.Item(30, 4).Value = -13.6
$ws.Cells.Item(30, 5).Value = -5.7
$ws.Cells.Item(30, 6).Value = 16.89
# Row 31
$ws.Cells.Item(31, 1).Value = "SC 132"
$ws.Cells.Item(31, 2).Value = -18.8
$ws.Cells.Item(31, 3).Value = 15.3
$ws.Cells.Item(31, 4).Value = -13.7
$ws.Cells.Item(31, 5).Value = -8.1
$ws.Cells.Item(31, 6).Value = 17.18
# Row 32
$ws.Cells.Item(32, 1).Value = "SC 193"
$ws.Cells.Item(32, 2).Value = -19.9
$ws.Cells.Item(32, 3).Value = 10.5
$ws.Cells.Item(32, 4).Value = -14.7
$ws.Cells.Item(32, 5).Value = -6.4
$ws.Cells.Item(32, 6).Value = 17.39
# Row 33
$ws.Cells.Item(33, 1).Value = "SC 232"
$ws.Cells.Item(33, 2).Value = -19.5
$ws.Cells.Item(33, 3).Value = 10.4
$ws.Cells.Item(33, 4).Value = -14.1
$ws.Cells.Item(33, 5).Value = -10.7
$ws.Cells.Item(33, 6).Value = 17.53
